# Insert a new weekly record row into the price list at row 93, pushing
# the existing rows 93:148 down to 94:149 (dimension grows from R148 to R149).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 93 (existing row 93 and below shift down by one).
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new weekly data point.
$ws.Cells.Item(93, 1).Value  = 6
$ws.Cells.Item(93, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(93, 3).Value  = "Metropolitana"
$ws.Cells.Item(93, 4).Value  = 44488
$ws.Cells.Item(93, 5).Value  = 13
$ws.Cells.Item(93, 6).Value  = 100112022
$ws.Cells.Item(93, 7).Value  = "Arveja Verde"
$ws.Cells.Item(93, 8).Value  = "Perfection"
$ws.Cells.Item(93, 9).Value  = "Primera"
$ws.Cells.Item(93, 10).Value = 180
$ws.Cells.Item(93, 11).Value = 20000
$ws.Cells.Item(93, 12).Value = 22000
$ws.Cells.Item(93, 13).Value = 20889
$ws.Cells.Item(93, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(93, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(93, 16).Value = 836
$ws.Cells.Item(93, 17).Value = 25
$ws.Cells.Item(93, 18).Value = "Hortaliza"
